$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "RUNMANAGER": add a new test-data row (row 11) for the new admin
# "add new user" test run.
# ---------------------------------------------------------------------------
$wsRun = $wb.Worksheets.Item("RUNMANAGER")
$wsRun.Range("A11").Value = "verifyThatTheAdminCanAddNewUser"
$wsRun.Range("B11").Value = "To check this test is executed"
$wsRun.Range("C11").Value = "yes"
$wsRun.Range("D11").Value = "'10"
$wsRun.Range("E11").Value = "'1"
[void]$wsRun.Range("F11").Select()

# ---------------------------------------------------------------------------
# Sheet "ADMINFUNCTIONALITY": new admin user, unique username/employee name,
# plus a new "newpassword" column for the left-menu change-password test.
# ---------------------------------------------------------------------------
$wsAdmin = $wb.Worksheets.Item("ADMINFUNCTIONALITY")

# Shift the old "confirmpassword" column from G to H and insert a new
# "newpassword" column in its place.
$wsAdmin.Range("H1").Value = "confirmpassword"
$wsAdmin.Range("H2").Value = "sham12345"
$wsAdmin.Range("H3").Value = "sham12345"
$wsAdmin.Range("G1").Value = "newpassword"
$wsAdmin.Range("G2").Value = "sham12345"
$wsAdmin.Range("G3").Value = "sham12345"

# Update username/employeename test data for the new admin user.
$wsAdmin.Range("C2").Value = "Admin"
$wsAdmin.Range("C3").Value = "Admin"
$wsAdmin.Range("D2").Value = "admin123"
$wsAdmin.Range("D3").Value = "admin123"
$wsAdmin.Range("E2").Value = "Orange  Test"
$wsAdmin.Range("E3").Value = "Orange  Test"

[void]$wsAdmin.Range("E3").Select()
